$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1531
$ws1.Range("F3").Value = 37
$ws1.Range("F6").Value = 2440
$ws1.Range("F8").Value = 1504
$ws1.Range("F10").Value = 176

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1531
$ws4.Range("F3").Value = 37
$ws4.Range("F6").Value = 2440
$ws4.Range("F8").Value = 1504
$ws4.Range("F10").Value = 176
